$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 21,20
$data[0,0] = "Sending cluster"
$data[0,1] = "Ligand symbol"
$data[0,2] = "Receptor symbol"
$data[0,3] = "Target cluster"
$data[0,4] = "Ligand-expressing cells"
$data[0,5] = "Ligand detection rate"
$data[0,6] = "Ligand average expression value"
$data[0,7] = "Ligand total expression value"
$data[0,8] = "Ligand derived specificity of average expression value"
$data[0,9] = "Ligand derived specificity of total expression value"
$data[0,10] = "Receptor-expressing cells"
$data[0,11] = "Receptor detection rate"
$data[0,12] = "Receptor average expression value"
$data[0,13] = "Receptor total expression value"
$data[0,14] = "Receptor derived specificity of average expression value"
$data[0,15] = "Receptor derived specificity of total expression value"
$data[0,16] = "Edge average expression weight"
$data[0,17] = "Edge total expression weight"
$data[0,18] = "Edge average expression derived specificity"
$data[0,19] = "Edge total expression derived specificity"
$data[1,0] = "ECs"
$data[1,1] = "Icosl"
$data[1,2] = "Icos"
$data[1,3] = "ECs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 0.5778633333333333
$data[1,7] = 1.73359
$data[1,8] = 0.01822222003587956
$data[1,9] = 0.01829736853065394
$data[1,10] = 1
$data[1,11] = 0.3333333333333333
$data[1,12] = 0.009971333333333334
$data[1,13] = 0.029914
$data[1,14] = 0.01215464113355011
$data[1,15] = 0.01292516977409686
$data[1,16] = 0.005762067917777777
$data[1,17] = 0.05185861126
$data[1,18] = 0.0002214845451927026
$data[1,19] = 0.0002364965946779193
$data[2,0] = "ECs"
$data[2,1] = "Icosl"
$data[2,2] = "Icos"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 0.5778633333333333
$data[2,7] = 1.73359
$data[2,8] = 0.01822222003587956
$data[2,9] = 0.01829736853065394
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.4560523333333333
$data[2,13] = 1.368157
$data[2,14] = 0.5559088503494856
$data[2,15] = 0.5911500134592177
$data[2,16] = 0.2635359215144444
$data[2,17] = 2.37182329363
$data[2,18] = 0.01012989339096117
$data[2,19] = 0.01081648965316434
$data[3,0] = "ECs"
$data[3,1] = "Icosl"
$data[3,2] = "Icos"
$data[3,3] = "MuSCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 0.5778633333333333
$data[3,7] = 1.73359
$data[3,8] = 0.01822222003587956
$data[3,9] = 0.01829736853065394
$data[3,10] = 2
$data[3,11] = 1
$data[3,12] = 0.1467185
$data[3,13] = 0.293437
$data[3,14] = 0.1788437569518725
$data[3,15] = 0.1267875591028168
$data[3,16] = 0.08478324147166666
$data[3,17] = 0.50869944883
$data[3,18] = 0.003258930291220385
$data[3,19] = 0.002319878694006306
$data[4,0] = "ECs"
$data[4,1] = "Icosl"
$data[4,2] = "Icos"
$data[4,3] = "Resolving-Mac"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 0.5778633333333333
$data[4,7] = 1.73359
$data[4,8] = 0.01822222003587956
$data[4,9] = 0.01829736853065394
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.2076303333333333
$data[4,13] = 0.622891
$data[4,14] = 0.2530927515650919
$data[4,15] = 0.2691372576638687
$data[4,16] = 0.1199819565211111
$data[4,17] = 1.07983760869
$data[4,18] = 0.004611911808505305
$data[4,19] = 0.004924503588805371
$data[5,0] = "FAPs"
$data[5,1] = "Icosl"
$data[5,2] = "Icos"
$data[5,3] = "ECs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.9308339999999999
$data[5,7] = 2.792502
$data[5,8] = 0.02935272232455987
$data[5,9] = 0.02947377304702275
$data[5,10] = 1
$data[5,11] = 0.3333333333333333
$data[5,12] = 0.009971333333333334
$data[5,13] = 0.029914
$data[5,14] = 0.01215464113355011
$data[5,15] = 0.01292516977409686
$data[5,16] = 0.009281656092000001
$data[5,17] = 0.083534904828
$data[5,18] = 0.0003567718061477699
$data[5,19] = 0.0003809535205159692
$data[6,0] = "FAPs"
$data[6,1] = "Icosl"
$data[6,2] = "Icos"
$data[6,3] = "Inflammatory-Mac"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.9308339999999999
$data[6,7] = 2.792502
$data[6,8] = 0.02935272232455987
$data[6,9] = 0.02947377304702275
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 0.4560523333333333
$data[6,13] = 1.368157
$data[6,14] = 0.5559088503494856
$data[6,15] = 0.5911500134592177
$data[6,16] = 0.4245090176459999
$data[6,17] = 3.820581158814
$data[6,18] = 0.01631743812207376
$data[6,19] = 0.01742342133344143
$data[7,0] = "FAPs"
$data[7,1] = "Icosl"
$data[7,2] = "Icos"
$data[7,3] = "MuSCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.9308339999999999
$data[7,7] = 2.792502
$data[7,8] = 0.02935272232455987
$data[7,9] = 0.02947377304702275
$data[7,10] = 2
$data[7,11] = 1
$data[7,12] = 0.1467185
$data[7,13] = 0.293437
$data[7,14] = 0.1788437569518725
$data[7,15] = 0.1267875591028168
$data[7,16] = 0.136570568229
$data[7,17] = 0.8194234093739999
$data[7,18] = 0.005249551137289387
$data[7,19] = 0.003736907742182405
$data[8,0] = "FAPs"
$data[8,1] = "Icosl"
$data[8,2] = "Icos"
$data[8,3] = "Resolving-Mac"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 0.9308339999999999
$data[8,7] = 2.792502
$data[8,8] = 0.02935272232455987
$data[8,9] = 0.02947377304702275
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.2076303333333333
$data[8,13] = 0.622891
$data[8,14] = 0.2530927515650919
$data[8,15] = 0.2691372576638687
$data[8,16] = 0.193269373698
$data[8,17] = 1.739424363282
$data[8,18] = 0.007428961259048957
$data[8,19] = 0.007932490450882951
$data[9,0] = "Inflammatory-Mac"
$data[9,1] = "Icosl"
$data[9,2] = "Icos"
$data[9,3] = "ECs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 16.208374
$data[9,7] = 48.625122
$data[9,8] = 0.5111114348580046
$data[9,9] = 0.5132192600799546
$data[9,10] = 1
$data[9,11] = 0.3333333333333333
$data[9,12] = 0.009971333333333334
$data[9,13] = 0.029914
$data[9,14] = 0.01215464113355011
$data[9,15] = 0.01292516977409686
$data[9,16] = 0.1616190999453334
$data[9,17] = 1.454571899508
$data[9,18] = 0.006212376069952918
$data[9,19] = 0.006633446067869785
$data[10,0] = "Inflammatory-Mac"
$data[10,1] = "Icosl"
$data[10,2] = "Icos"
$data[10,3] = "Inflammatory-Mac"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 16.208374
$data[10,7] = 48.625122
$data[10,8] = 0.5111114348580046
$data[10,9] = 0.5132192600799546
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 0.4560523333333333
$data[10,13] = 1.368157
$data[10,14] = 0.5559088503494856
$data[10,15] = 0.5911500134592177
$data[10,16] = 7.391866782239334
$data[10,17] = 66.526801040154
$data[10,18] = 0.2841313701523893
$data[10,19] = 0.3033895725037949
$data[11,0] = "Inflammatory-Mac"
$data[11,1] = "Icosl"
$data[11,2] = "Icos"
$data[11,3] = "MuSCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 16.208374
$data[11,7] = 48.625122
$data[11,8] = 0.5111114348580046
$data[11,9] = 0.5132192600799546
$data[11,10] = 2
$data[11,11] = 1
$data[11,12] = 0.1467185
$data[11,13] = 0.293437
$data[11,14] = 0.1788437569518725
$data[11,15] = 0.1267875591028168
$data[11,16] = 2.378068320719001
$data[11,17] = 14.268409924314
$data[11,18] = 0.09140908923106777
$data[11,19] = 0.06506981727009113
$data[12,0] = "Inflammatory-Mac"
$data[12,1] = "Icosl"
$data[12,2] = "Icos"
$data[12,3] = "Resolving-Mac"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 16.208374
$data[12,7] = 48.625122
$data[12,8] = 0.5111114348580046
$data[12,9] = 0.5132192600799546
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 0.2076303333333333
$data[12,13] = 0.622891
$data[12,14] = 0.2530927515650919
$data[12,15] = 0.2691372576638687
$data[12,16] = 3.365350096411334
$data[12,17] = 30.288150867702
$data[12,18] = 0.1293585994045946
$data[12,19] = 0.1381264242381988
$data[13,0] = "MuSCs"
$data[13,1] = "Icosl"
$data[13,2] = "Icos"
$data[13,3] = "ECs"
$data[13,4] = 2
$data[13,5] = 1
$data[13,6] = 0.39073
$data[13,7] = 0.78146
$data[13,8] = 0.01232119711342224
$data[13,9] = 0.008248006513630573
$data[13,10] = 1
$data[13,11] = 0.3333333333333333
$data[13,12] = 0.009971333333333334
$data[13,13] = 0.029914
$data[13,14] = 0.01215464113355011
$data[13,15] = 0.01292516977409686
$data[13,16] = 0.003896099073333334
$data[13,17] = 0.02337659444
$data[13,18] = 0.0001497597292493808
$data[13,19] = 0.0001066068844865319
$data[14,0] = "MuSCs"
$data[14,1] = "Icosl"
$data[14,2] = "Icos"
$data[14,3] = "Inflammatory-Mac"
$data[14,4] = 2
$data[14,5] = 1
$data[14,6] = 0.39073
$data[14,7] = 0.78146
$data[14,8] = 0.01232119711342224
$data[14,9] = 0.008248006513630573
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 0.4560523333333333
$data[14,13] = 1.368157
$data[14,14] = 0.5559088503494856
$data[14,15] = 0.5911500134592177
$data[14,16] = 0.1781933282033333
$data[14,17] = 1.06915996922
$data[14,18] = 0.006849462522251959
$data[14,19] = 0.004875809161544428
$data[15,0] = "MuSCs"
$data[15,1] = "Icosl"
$data[15,2] = "Icos"
$data[15,3] = "MuSCs"
$data[15,4] = 2
$data[15,5] = 1
$data[15,6] = 0.39073
$data[15,7] = 0.78146
$data[15,8] = 0.01232119711342224
$data[15,9] = 0.008248006513630573
$data[15,10] = 2
$data[15,11] = 1
$data[15,12] = 0.1467185
$data[15,13] = 0.293437
$data[15,14] = 0.1788437569518725
$data[15,15] = 0.1267875591028168
$data[15,16] = 0.057327319505
$data[15,17] = 0.22930927802
$data[15,18] = 0.002203569181909
$data[15,19] = 0.001045744613327354
$data[16,0] = "MuSCs"
$data[16,1] = "Icosl"
$data[16,2] = "Icos"
$data[16,3] = "Resolving-Mac"
$data[16,4] = 2
$data[16,5] = 1
$data[16,6] = 0.39073
$data[16,7] = 0.78146
$data[16,8] = 0.01232119711342224
$data[16,9] = 0.008248006513630573
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 0.2076303333333333
$data[16,13] = 0.622891
$data[16,14] = 0.2530927515650919
$data[16,15] = 0.2691372576638687
$data[16,16] = 0.08112740014333333
$data[16,17] = 0.48676440086
$data[16,18] = 0.003118405680011903
$data[16,19] = 0.002219845854272259
$data[17,0] = "Resolving-Mac"
$data[17,1] = "Icosl"
$data[17,2] = "Icos"
$data[17,3] = "ECs"
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 13.60421466666667
$data[17,7] = 40.81264400000001
$data[17,8] = 0.4289924256681337
$data[17,9] = 0.4307615918287382
$data[17,10] = 1
$data[17,11] = 0.3333333333333333
$data[17,12] = 0.009971333333333334
$data[17,13] = 0.029914
$data[17,14] = 0.01215464113355011
$data[17,15] = 0.01292516977409686
$data[17,16] = 0.1356521591795556
$data[17,17] = 1.220869432616
$data[17,18] = 0.005214248983007335
$data[17,19] = 0.005567666706546656
$data[18,0] = "Resolving-Mac"
$data[18,1] = "Icosl"
$data[18,2] = "Icos"
$data[18,3] = "Inflammatory-Mac"
$data[18,4] = 3
$data[18,5] = 1
$data[18,6] = 13.60421466666667
$data[18,7] = 40.81264400000001
$data[18,8] = 0.4289924256681337
$data[18,9] = 0.4307615918287382
$data[18,10] = 3
$data[18,11] = 1
$data[18,12] = 0.4560523333333333
$data[18,13] = 1.368157
$data[18,14] = 0.5559088503494856
$data[18,15] = 0.5911500134592177
$data[18,16] = 6.20423384190089
$data[18,17] = 55.838104577108
$data[18,18] = 0.2384806861618094
$data[18,19] = 0.2546447208072726
$data[19,0] = "Resolving-Mac"
$data[19,1] = "Icosl"
$data[19,2] = "Icos"
$data[19,3] = "MuSCs"
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 13.60421466666667
$data[19,7] = 40.81264400000001
$data[19,8] = 0.4289924256681337
$data[19,9] = 0.4307615918287382
$data[19,10] = 2
$data[19,11] = 1
$data[19,12] = 0.1467185
$data[19,13] = 0.293437
$data[19,14] = 0.1788437569518725
$data[19,15] = 0.1267875591028168
$data[19,16] = 1.995989969571334
$data[19,17] = 11.975939817428
$data[19,18] = 0.07672261711038593
$data[19,19] = 0.05461521078320958
$data[20,0] = "Resolving-Mac"
$data[20,1] = "Icosl"
$data[20,2] = "Icos"
$data[20,3] = "Resolving-Mac"
$data[20,4] = 3
$data[20,5] = 1
$data[20,6] = 13.60421466666667
$data[20,7] = 40.81264400000001
$data[20,8] = 0.4289924256681337
$data[20,9] = 0.4307615918287382
$data[20,10] = 3
$data[20,11] = 1
$data[20,12] = 0.2076303333333333
$data[20,13] = 0.622891
$data[20,14] = 0.2530927515650919
$data[20,15] = 0.2691372576638687
$data[20,16] = 2.824647625978223
$data[20,17] = 25.421828633804
$data[20,18] = 0.1085748734129311
$data[20,19] = 0.1159339935317093

$ws.Range("A1:T21").Value = $data
